$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.186.04"
$ws.Range("D3").Value = "3.703.95"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.20"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.91"
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "661.80"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.426"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("D11").Value = "3.699.27"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000319"
$ws.Range("E12").Value = "  +19.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "44.26"
$ws.Range("E13").Value = "  -3.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.208"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.81"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "4.390.70"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").Value = "96.838.98"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.06"
$ws.Range("E18").Value = "  +1.88%  "
$ws.Range("D19").Value = "3.682.00"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.95"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.68"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.505"
$ws.Range("E22").Value = "  -3.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "521.77"
$ws.Range("E23").Value = "  -1.01%  "
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000219"
$ws.Range("E25").Value = "  +7.27%  "
$ws.Range("E26").Value = "  -2.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.89"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.194"
$ws.Range("E28").Value = "  +15.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.55"
$ws.Range("E29").Value = "  +3.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.80"
$ws.Range("E30").Value = "  +2.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.03"
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.28%  "
$ws.Range("E33").Value = "  +2.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.86"
$ws.Range("E34").Value = "  -3.61%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "654.24"
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "32.21"
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.593"
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.87"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.166"
$ws.Range("E41").Value = "  +3.50%  "
$ws.Range("E42").Value = "  +5.01%  "
$ws.Range("E43").Value = "  +3.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.18"
$ws.Range("E44").Value = "  -9.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.480"
$ws.Range("E45").Value = "  +7.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.969"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0457"
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.32"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.62"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.68"
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.55"
$ws.Range("E51").Value = "  +1.03%  "
